$d = $word.ActiveDocument

$pairs = @(
    @("28×31=", "89×21="),
    @("52×43=", "78×75="),
    @("12×89=", "60×26="),
    @("48×41=", "91×51="),
    @("15×70=", "39×76="),
    @("54×95=", "94×49="),
    @("39×46=", "82×49="),
    @("85×36=", "76×42="),
    @("83×89=", "30×63="),
    @("13×27=", "57×50="),
    @("86×42=", "88×60="),
    @("63×76=", "57×19="),
    @("43×82=", "74×91="),
    @("81×89=", "65×27="),
    @("40×54=", "58×45="),
    @("14×44=", "78×81="),
    @("30×40=", "66×16="),
    @("15×95=", "46×37="),
    @("72×97=", "12×77="),
    @("39×57=", "62×99="),
    @("42×52=", "60×57="),
    @("38×88=", "11×60="),
    @("21×69=", "41×99="),
    @("11×38=", "40×94="),
    @("39×89=", "83×50=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
